# edit.ps1
# Applies the "next day" stock-report delta described by the commit diff:
#   - A handful of quantity-on-hand (column F) corrections (mostly -1, a few
#     larger swings), with a few item-code (column B) / rate (D,E) corrections
#     where two rows for the same item label had their records swapped.
#   - Column G (Value) is always Rate(D) * Qty(F), so it is recomputed for
#     every touched row rather than hard-coded.
#   - Each company section's "Sub Total:" row (column B) is the sum of G over
#     that section's data rows, so it is recomputed from the refreshed data.
#   - The final "Sub Total:" row (938) mirrors the "Grand Total:" row (939),
#     which is the sum of all 80 section subtotals; both are recomputed too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Direct cell corrections: @(row, col, newValue) ------------------
# col 2=B (Item Code), 4=D (Rate), 5=E (MRP), 6=F (Qty)
$cellEdits = @(
    @(30, 6, 69),
    @(32, 2, 51755),
    @(32, 6, 1),
    @(33, 2, 66452),
    @(33, 6, 68),
    @(38, 6, 44),
    @(97, 6, 16),
    @(130, 6, 45),
    @(135, 6, 111),
    @(186, 6, 57),
    @(222, 6, 85),
    @(234, 2, 57004),
    @(234, 6, 5),
    @(235, 2, 63255),
    @(235, 6, 83),
    @(276, 6, 246),
    @(279, 6, 111),
    @(284, 6, 51),
    @(285, 6, 2027),
    @(290, 6, 94),
    @(310, 6, 11),
    @(314, 6, 10),
    @(320, 6, 83),
    @(322, 6, 7),
    @(323, 6, 42),
    @(372, 6, 6),
    @(406, 6, 182),
    @(457, 6, 366),
    @(466, 6, 40),
    @(469, 6, 312),
    @(481, 2, 58047),
    @(481, 4, 105.54),
    @(481, 5, 126.1),
    @(481, 6, 34),
    @(482, 2, 47097),
    @(482, 4, 112.28),
    @(482, 5, 134.16),
    @(482, 6, 15),
    @(487, 6, 576),
    @(509, 6, 272),
    @(519, 6, 94),
    @(560, 2, 64927),
    @(560, 5, 17.26),
    @(560, 6, 106),
    @(561, 2, 45718),
    @(561, 5, 19.38),
    @(561, 6, -294),
    @(584, 6, 18),
    @(620, 6, 44),
    @(641, 2, 64810),
    @(641, 5, 291.22),
    @(641, 6, 2),
    @(642, 2, 53319),
    @(642, 5, 310.64),
    @(642, 6, -6),
    @(677, 6, 236),
    @(678, 6, 354),
    @(680, 6, 234),
    @(682, 6, 13),
    @(684, 6, 152),
    @(686, 6, 513),
    @(701, 6, 56),
    @(703, 6, 53),
    @(716, 6, 9),
    @(774, 6, 0),
    @(805, 6, 116),
    @(810, 6, 16),
    @(823, 6, 442),
    @(825, 6, 556),
    @(826, 6, 99),
    @(827, 2, 65079),
    @(827, 6, 6),
    @(828, 2, 65362),
    @(828, 6, 0),
    @(829, 6, 401),
    @(848, 6, 1),
    @(875, 6, 91),
    @(880, 6, 11),
    @(887, 6, 1604),
    @(903, 6, 26)
)

foreach ($edit in $cellEdits) {
    $row = $edit[0]
    $col = $edit[1]
    $val = $edit[2]
    $ws.Cells.Item($row, $col).Value = $val
}

# --- 2. Recompute Value (G) = Rate (D) * Qty (F) for every touched row --
$affectedRows = @(30, 32, 33, 38, 97, 130, 135, 186, 222, 234, 235, 276, 279, 284, 285, 290, 310, 314, 320, 322, 323, 372, 406, 457, 466, 469, 481, 482, 487, 509, 519, 560, 561, 584, 620, 641, 642, 677, 678, 680, 682, 684, 686, 701, 703, 716, 774, 805, 810, 823, 825, 826, 827, 828, 829, 848, 875, 880, 887, 903)

foreach ($row in $affectedRows) {
    $rate = $ws.Cells.Item($row, 4).Value2
    $qty  = $ws.Cells.Item($row, 6).Value2
    $ws.Cells.Item($row, 7).Value = $rate * $qty
}

# --- 3. Recompute each affected company "Sub Total:" (col B) ------------
# @(subtotalRow, dataStartRow, dataEndRow)
$subtotalRanges = @(
    @(40, 21, 39),
    @(101, 78, 100),
    @(146, 117, 145),
    @(190, 184, 189),
    @(248, 220, 247),
    @(281, 273, 280),
    @(292, 283, 291),
    @(316, 306, 315),
    @(373, 318, 372),
    @(410, 375, 409),
    @(459, 454, 458),
    @(471, 461, 470),
    @(490, 479, 489),
    @(520, 492, 519),
    @(585, 574, 584),
    @(634, 620, 633),
    @(693, 677, 692),
    @(706, 698, 705),
    @(730, 715, 729),
    @(784, 770, 783),
    @(836, 802, 835),
    @(870, 845, 869),
    @(881, 872, 880),
    @(893, 886, 892),
    @(909, 901, 908)
)

foreach ($range in $subtotalRanges) {
    $subtotalRow = $range[0]
    $dataStart = $range[1]
    $dataEnd = $range[2]
    $sum = 0.0
    for ($r = $dataStart; $r -le $dataEnd; $r++) {
        $sum = $sum + $ws.Cells.Item($r, 7).Value2
    }
    $ws.Cells.Item($subtotalRow, 2).Value = $sum
}

# --- 4. Recompute Grand Total (row 939) as sum of all 80 section ---------
#        subtotals, and mirror it into the duplicate "Sub Total:" row 938.
$allSubtotalRows = @(10, 13, 19, 40, 73, 76, 101, 105, 108, 115, 146, 151, 159, 162, 172, 182, 190, 200, 218, 248, 253, 257, 262, 271, 281, 292, 304, 316, 373, 410, 414, 422, 428, 452, 459, 471, 474, 477, 490, 520, 532, 546, 572, 585, 588, 592, 596, 601, 613, 618, 634, 637, 644, 657, 672, 675, 693, 696, 706, 709, 713, 730, 740, 753, 762, 765, 768, 784, 793, 800, 836, 843, 870, 881, 884, 893, 899, 909, 932, 937)

$grandTotal = 0.0
foreach ($r in $allSubtotalRows) {
    $grandTotal = $grandTotal + $ws.Cells.Item($r, 2).Value2
}

$ws.Cells.Item(938, 2).Value = $grandTotal
$ws.Cells.Item(939, 2).Value = $grandTotal
